{"js": "// Helper: search (optionally scoped to a given range/body) for `searchText`\n// and replace the FIRST match with `replaceText`. Throws if not found so\n// problems surface loudly instead of silently no-op'ing.\nasync function replaceFirst(scope, searchText, replaceText, options) {\n  const results = scope.search(searchText, Object.assign({ matchCase: true, matchWholeWord: false }, options || {}));\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(replaceText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Helper: find the first paragraph whose text contains `anchorText` (used to\n// scope subsequent, possibly-ambiguous, searches to the right paragraph).\nasync function paragraphContaining(body, anchorText) {\n  const results = body.search(anchorText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Anchor text not found: \" + anchorText);\n  }\n  const para = results.items[0].paragraphs.getFirst();\n  return para;\n}\n\nconst body = context.document.body;\n\n// 1. \"...will almost certainly scanned...\" -> \"...will almost certainly be scanned...\"\nawait replaceFirst(\n  body,\n  \"will almost certainly scanned\",\n  \"will almost certainly be scanned\"\n);\n\n// 2. \"...which they keep secure at all times...\" -> \"...which they always keep secure...\"\nawait replaceFirst(\n  body,\n  \"which they keep secure at all times\",\n  \"which they always keep secure\"\n);\n\n// 3. \"... RSA for SSH v2 (or use -t ecdsa or -t ed25519.) ...\" -> \"... (or -t ed25519.) ...\"\n//    Done as three small, run-boundary-respecting edits so the Courier New\n//    formatting of the command text survives the edit.\n{\n  const para = await paragraphContaining(body, \"Only create the first key\");\n\n  await replaceFirst(para, \" (or use \", \" (or \");\n\n  // Remove the word \"ecdsa\" (its own run/proofErr span).\n  const ecdsaResults = para.search(\"ecdsa\", { matchCase: true });\n  ecdsaResults.load(\"items\");\n  await context.sync();\n  if (ecdsaResults.items.length === 0) {\n    throw new Error(\"Text not found: ecdsa (in RSA key paragraph)\");\n  }\n  ecdsaResults.items[0].delete();\n  await context.sync();\n\n  // Collapse \" or -t ed25519\" down to \"ed25519\" (keeps the leading \"-t \" run).\n  await replaceFirst(para, \" or -t ed25519\", \"ed25519\");\n}\n\n// 4. \"...security than rsa providing...\" -> \"...security than RSA, providing...\"\nawait replaceFirst(\n  body,\n  \"security than rsa providing\",\n  \"security than RSA, providing\"\n);\n\n// 5. \"...SSH keys in Linix is...\" -> \"...SSH keys in Linux is...\"\nawait replaceFirst(\n  body,\n  \"SSH keys in Linix is\",\n  \"SSH keys in Linux is\"\n);\n\n// 6. \"The most popular SSH Client...\" -> \"It used to be that the most popular SSH Client...\"\nawait replaceFirst(\n  body,\n  \"The most popular SSH Client\",\n  \"It used to be that the most popular SSH Client\"\n);\n\n// 7. Add a new \"Hand In\" section at the end of the document.\nconst handInHeading = body.insertParagraph(\"Hand In\", Word.InsertLocation.end);\nhandInHeading.styleBuiltIn = Word.Style.heading1;\nawait context.sync();\n\nconst handInBody = body.insertParagraph(\n  \"Submit a screenshot of you connecting to your SSH server using private/public keys instead of a password.\",\n  Word.InsertLocation.end\n);\nhandInBody.styleBuiltIn = Word.Style.normal;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Replace the FIRST occurrence of $findText with $replaceText. wdReplaceOne (1)\n# is used (not wdReplaceAll) so ambiguous search strings only touch the single\n# intended occurrence.\nfunction Replace-FirstText {\n    param(\n        [string]$findText,\n        [string]$replaceText\n    )\n    $rng = $d.Content\n    $ok = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n    if (-not $ok) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n# 1. \"...will almost certainly scanned...\" -> \"...will almost certainly be scanned...\"\nReplace-FirstText \"will almost certainly scanned\" \"will almost certainly be scanned\"\n\n# 2. \"...which they keep secure at all times...\" -> \"...which they always keep secure...\"\nReplace-FirstText \"which they keep secure at all times\" \"which they always keep secure\"\n\n# 3. \"... RSA for SSH v2 (or use -t ecdsa or -t ed25519.) ...\" -> \"... (or -t ed25519.) ...\"\n#    Done as three small edits (rather than one big replace) so the Courier\n#    New formatting of the command text survives the edit.\nReplace-FirstText \" (or use \" \" (or \"\nReplace-FirstText \"ecdsa\" \"\"\nReplace-FirstText \" or -t ed25519\" \"ed25519\"\n\n# 4. \"...security than rsa providing...\" -> \"...security than RSA, providing...\"\nReplace-FirstText \"security than rsa providing\" \"security than RSA, providing\"\n\n# 5. \"...SSH keys in Linix is...\" -> \"...SSH keys in Linux is...\"\nReplace-FirstText \"SSH keys in Linix is\" \"SSH keys in Linux is\"\n\n# 6. \"The most popular SSH Client...\" -> \"It used to be that the most popular SSH Client...\"\nReplace-FirstText \"The most popular SSH Client\" \"It used to be that the most popular SSH Client\"\n\n# 7. Add a new \"Hand In\" section at the end of the document.\n$endRange = $d.Content\n$endRange.Collapse(0)\n$endRange.InsertParagraphAfter()\n$endRange.Collapse(0)\n$endRange.Text = \"Hand In\"\n$d.Paragraphs.Last.Style = \"Heading 1\"\n\n$endRange2 = $d.Content\n$endRange2.Collapse(0)\n$endRange2.InsertParagraphAfter()\n$endRange2.Collapse(0)\n$endRange2.Text = \"Submit a screenshot of you connecting to your SSH server using private/public keys instead of a password.\"\n$d.Paragraphs.Last.Style = \"Normal\"\n"}
